$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Column B (dietas): now holds the full alphabetically sorted diet list ---
$ws.Range("B2").Value = 'dairy free '
$ws.Range("B3").Value = 'gluten free '
$ws.Range("B4").Value = 'ketogenic '
$ws.Range("B5").Value = 'lacto ovo vegetarian '
$ws.Range("B6").Value = 'lacto vegetarian '
$ws.Range("B7").Value = 'low fodmap '
$ws.Range("B8").Value = 'ovo vegetarian '
$ws.Range("B9").Value = 'paleolithic '
$ws.Range("B10").Value = 'pescetarian '
$ws.Range("B11").Value = 'primal '
$ws.Range("B12").Value = 'vegan '
$ws.Range("B13").Value = 'vegetarian '
$ws.Range("B14").Value = 'whole 30'

# --- G2: updated steps for "Pizza a la francesa" ---
$ws.Range("G2").Value = '{"Comprar 50 gr de masa para pizza lista","extender en molde enmantecado","cubrir con salsa para pizza al gusto","agregar 150 gr de queso mozzarella","incorporar 50 gr de queso Roquefort","agregar los ongos y las especia al gusto","hornear a 180 por 15 min"}'

# --- H8: corrected image URL (drop stray "!d" suffix) ---
$ws.Range("H8").Value = 'https://c.wallhere.com/photos/de/ca/food_sandwiches-216633.jpg'

# --- New hyperlinks for H3:H8, mirroring the existing H2 pattern ---
$ws.Hyperlinks.Add($ws.Range("H3"), 'https://www.cocinayvino.com/wp-content/uploads/2017/09/MANDOCA-e1504641550954.jpg') | Out-Null
$ws.Range("H3").Style = $ws.Range("H2").Style
$ws.Hyperlinks.Add($ws.Range("H4"), 'https://d320djwtwnl5uo.cloudfront.net/recetas/cover/cooki_lg91oLQjnwTD8XNvRqYI2MEdO34xct.png') | Out-Null
$ws.Range("H4").Style = $ws.Range("H2").Style
$ws.Hyperlinks.Add($ws.Range("H5"), 'https://placeralplato.com/files/2015/09/Galletas-de-nuez.jpg') | Out-Null
$ws.Range("H5").Style = $ws.Range("H2").Style
$ws.Hyperlinks.Add($ws.Range("H6"), 'https://t2.rg.ltmcdn.com/es/posts/6/2/9/galletas_con_chispas_de_chocolate_caseras_35926_600.jpg') | Out-Null
$ws.Range("H6").Style = $ws.Range("H2").Style
$ws.Hyperlinks.Add($ws.Range("H7"), 'https://t2.uc.ltmcdn.com/es/posts/4/3/1/como_hacer_arepas_venezolanas_28134_600.jpg') | Out-Null
$ws.Range("H7").Style = $ws.Range("H2").Style
$ws.Hyperlinks.Add($ws.Range("H8"), 'https://c.wallhere.com/photos/de/ca/food_sandwiches-216633.jpg') | Out-Null
$ws.Range("H8").Style = $ws.Range("H2").Style

# --- B4 gets a light-gray fill to highlight it ---
$ws.Range("B4").Interior.Color = 10921638

# --- Selection ends on H7, matching the authored workbook state ---
$ws.Activate()
$ws.Range("H7").Select()
